$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (dates are sequential days after the last existing date)
$data = @(
    @(45630, 200),
    @(45631, 500),
    @(45632, 50),
    @(45633, 350),
    @(45634, 100),
    @(45635, 175),
    @(45636, 350)
)

$row = 5
foreach ($pair in $data) {
    # Copy the formatting (borders, wrap, style) from the row above so new rows
    # match the existing table styling.
    $ws.Range("A$($row - 1):B$($row - 1)").Copy()
    $ws.Range("A$($row):B$($row)").PasteSpecial(-4122)  # xlPasteFormats

    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]

    $row++
}

$excel.CutCopyMode = 0

# Update the active cell selection to reflect the cell after the new last row
$ws.Range("B12").Select()
